# Adds new LeetCode practice rows (16-27) and fixes the Lowest Common
# Ancestor hyperlink (row 10) - "added problems till kclosest points to origin"

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 10: fix the URL (add trailing slash) and wire it up as a real
#     hyperlink, matching the style used by the other hyperlink cells in
#     column E (Hiperlink cell style + vertical-center alignment). ---
$ws.Range("E10").Value = "https://leetcode.com/problems/lowest-common-ancestor-of-a-binary-tree/"
$ws.Hyperlinks.Add($ws.Range("E10"), "https://leetcode.com/problems/lowest-common-ancestor-of-a-binary-tree/")
$ws.Range("E10").Style = "Hiperlink"
$ws.Range("E10").VerticalAlignment = -4108

# --- New rows of practice problems ---
$rows = @(
  @{ r=16; a="Ransom Note";                    b="string";      c="counter";              d="";    e="https://leetcode.com/problems/ransom-note/" },
  @{ r=17; a="Climbing Stairs";                 b="recursive";   c="dynamic programing";   d="";    e="https://leetcode.com/problems/climbing-stairs/" },
  @{ r=18; a="Longest palindrome";               b="string";      c="counter";              d="";    e="https://leetcode.com/problems/longest-palindrome/" },
  @{ r=19; a="Majority Element";                 b="array";       c="counter";              d="";    e="https://leetcode.com/problems/majority-element/" },
  @{ r=20; a="Add Binary";                       b="binary math"; c="";                     d="";    e="https://leetcode.com/problems/add-binary/" },
  @{ r=21; a="Diameter of Binary Tree";           b="tree";        c="recursive";            d="dfs"; e="https://leetcode.com/problems/diameter-of-binary-tree/" },
  @{ r=22; a="Middle of the Linked List";         b="linked list"; c="two pointers";         d="";    e="https://leetcode.com/problems/middle-of-the-linked-list/" },
  @{ r=23; a="Maximum Depth of Binary Tree";      b="tree";        c="dfs";                  d="";    e="https://leetcode.com/problems/maximum-depth-of-binary-tree/" },
  @{ r=24; a="Contains Duplicate";                b="array";       c="";                     d="";    e="https://leetcode.com/problems/contains-duplicate/" },
  @{ r=25; a="Maximum Subarray";                  b="array";       c="";                     d="";    e="https://leetcode.com/problems/maximum-subarray/" },
  @{ r=26; a="Insert Interval";                   b="array";       c="";                     d="";    e="https://leetcode.com/problems/insert-interval/" },
  @{ r=27; a="01 matrix";                         b="matrix";      c="stack";                d="bfs"; e="https://leetcode.com/problems/01-matrix/" }
)

foreach ($row in $rows) {
  $r = $row.r
  $ws.Cells.Item($r, 1).Value = $row.a
  $ws.Cells.Item($r, 2).Value = $row.b
  if ($row.c -ne "") { $ws.Cells.Item($r, 3).Value = $row.c }
  if ($row.d -ne "") { $ws.Cells.Item($r, 4).Value = $row.d }
  $ws.Cells.Item($r, 5).Value = $row.e
}

# --- View state: move selection to just past the new data, matching the
#     author's cursor position after pasting the new rows. ---
$ws.Range("C28").Select()

Write-Host "applied"
